{"js": "const doc = context.document;\nconst body = doc.body;\n\n// ---------------------------------------------------------------------\n// Change 1: Move the \"_GoBack\" bookmark to sit right after the title\n// text \"WANGRY:  Are We Angry For Water\" (an empty/collapsed bookmark).\n// The document only ever has a single \"_GoBack\" bookmark (it marks the\n// last edit location), so we first drop the one currently sitting\n// inside \"JT.ipynb\" further down, then insert a fresh one in the new\n// spot.\n// ---------------------------------------------------------------------\ndoc.deleteBookmark(\"_GoBack\");\n\nconst titleResults = body.search(\"WANGRY:  Are We Angry For Water\", { matchCase: true });\nawait context.sync();\nconst titleRange = titleResults.items[0];\nconst afterTitle = titleRange.getRange(\"After\");\nafterTitle.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Change 2: Extend the \"All notebooks ... Load, Manipulation, Out\"\n// sentence with a new trailing sentence, as its own run. Plain\n// insertText onto a collapsed range just grows the existing run's text\n// (same formatting), so re-locate the freshly inserted text via a new\n// search() and nudge Bold on/off to force the engine to split it into\n// a distinct w:r (the run ends up with the very same rPr either way,\n// since Bold is turned back off again).\n// ---------------------------------------------------------------------\nconst sentenceResults = body.search(\n  \"All notebooks will have in addition what it is doing:  Load, Manipulation, Out\",\n  { matchCase: true }\n);\nawait context.sync();\nconst sentenceRange = sentenceResults.items[0];\nconst afterSentence = sentenceRange.getRange(\"After\");\nconst newSentence = \".  You can figure it out based on the folder you will store the code in.\";\nafterSentence.insertText(newSentence, \"Start\");\nawait context.sync();\n\nconst newRunResults = body.search(newSentence, { matchCase: true });\nawait context.sync();\nconst newRunRange = newRunResults.items[0];\nnewRunRange.font.bold = true;\nnewRunRange.font.bold = false;\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Change 3: Merge the \"JT.\" + \"i\" + \"pynb\" runs (previously split\n// around the old \"_GoBack\" bookmark) into a single run reading\n// \"JT.ipynb\". Replacing the whole matched range's text in one shot\n// keeps the surrounding w:proofErr markers in their original positions\n// and keeps the bold/size run formatting intact.\n// ---------------------------------------------------------------------\nconst fileNameResults = body.search(\"JT.ipynb\", { matchCase: true });\nawait context.sync();\nconst fileNameRange = fileNameResults.items[0];\nfileNameRange.insertText(\"JT.ipynb\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Change 1: Move the \"_GoBack\" bookmark to sit right after the title\n# text \"WANGRY:  Are We Angry For Water\" (empty/collapsed bookmark).\n# Word only ever keeps a single \"_GoBack\" bookmark, so adding a new one\n# with that name automatically drops whichever one previously existed\n# (the one currently sitting inside \"JT.ipynb\" further down).\n#\n# A genuinely collapsed Range sitting exactly at \"end of paragraph text,\n# right before the paragraph mark\" can't be used directly with\n# Bookmarks.Add, so we temporarily insert a one-character placeholder,\n# wrap the bookmark around that (non-collapsed) placeholder range, and\n# then delete just the placeholder character again - leaving the empty\n# bookmark pair behind in the right spot.\n# ---------------------------------------------------------------------\n$titlePara = $d.Paragraphs.Item(1)\n$insertPos = $titlePara.Range.End - 1\n$placeholder = $d.Range($insertPos, $insertPos)\n$placeholder.InsertAfter(\"X\")\n$bookmarkRange = $d.Range($insertPos, $insertPos + 1)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n$d.Range($insertPos, $insertPos + 1).Delete()\n\n# ---------------------------------------------------------------------\n# Change 2: Extend the \"All notebooks ... Load, Manipulation, Out\"\n# sentence with a new trailing sentence, as its own run (matching\n# formatting). InsertAfter on a collapsed range with identical run\n# formatting just extends the existing run's text, so we nudge the\n# Bold flag on/off on the freshly inserted span to force Word to keep\n# it as a distinct w:r (the rPr ends up identical either way, since\n# Bold is restored to its original/absent state).\n# ---------------------------------------------------------------------\n$target = $d.Content\n$target.Find.Execute(\"All notebooks will have in addition what it is doing:  Load, Manipulation, Out\") | Out-Null\n$target.Collapse(0)\n$newSentence = \".  You can figure it out based on the folder you will store the code in.\"\n$insertStart = $target.Start\n$target.InsertAfter($newSentence)\n$newRunRange = $d.Range($insertStart, $insertStart + $newSentence.Length)\n$newRunRange.Font.Bold = 1\n$newRunRange.Font.Bold = 0\n\n# ---------------------------------------------------------------------\n# Change 3: Merge the \"JT.\" + \"i\" + \"pynb\" runs (previously split around\n# the old \"_GoBack\" bookmark) into a single run reading \"JT.ipynb\".\n# Using Find/Replace (rather than deleting+reinserting) keeps the\n# surrounding w:proofErr markers in their original positions and keeps\n# the bold/size run formatting intact.\n# ---------------------------------------------------------------------\n$target2 = $d.Content\n$target2.Find.Execute(\"JT.ipynb\", $false, $false, $false, $false, $false, $true, 1, $false, \"JT.ipynb\", 2) | Out-Null\n"}
